# Insert a new data row at row 79 (pushing the existing rows 79-145 down to
# 80-146) and populate the new row with the new "Haba" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(79).Insert()

$ws.Cells.Item(79, 1).Value = 5
$ws.Cells.Item(79, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(79, 3).Value = "Maule"
$ws.Cells.Item(79, 4).Value = 45240
$ws.Cells.Item(79, 5).Value = 7
$ws.Cells.Item(79, 6).Value = 100112026
$ws.Cells.Item(79, 7).Value = "Haba"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 400
$ws.Cells.Item(79, 11).Value = 9000
$ws.Cells.Item(79, 12).Value = 9000
$ws.Cells.Item(79, 13).Value = 9000
$ws.Cells.Item(79, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(79, 15).Value = "Región del Maule"
$ws.Cells.Item(79, 16).Value = 360
$ws.Cells.Item(79, 17).Value = 25
$ws.Cells.Item(79, 18).Value = "Hortaliza"
